$d = $word.ActiveDocument

# --- Locate the target paragraph: the Source Code paragraph that begins
#     with ". pt_base age gender ethnicity, post(`postname')" ---------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("pt_base age gender ethnicity", $true, $false, $false, `
                        $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target paragraph for edit"
}

$matchRange = $find.Parent
$targetPara = $matchRange.Paragraphs(1)
$insertPos  = $targetPara.Range.Start

# --- Insert the new "post" example line, styled like the existing
#     Verbatim/SourceCode text, followed by a manual line break, right
#     before the existing ". pt_base ..." run. ------------------------------
$newLineText = "      . post ``postname' (`"Variable`") (`"Summary`")"

$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertBefore($newLineText)
$insertRange.Style = $d.Styles("Verbatim Char")

$breakPos = $insertRange.End
$breakRange = $d.Range($breakPos, $breakPos)
$breakRange.InsertBreak(6)
